$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "范总"
$ws.Range("B6").Value = "王柳"
$ws.Range("C6").Value = "张三"
$ws.Range("D6").Value = "王五"
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 12
$ws.Range("G6").Value = 66

$ws.Range("A7").Value = "杨总"
$ws.Range("B7").Value = "李四"
$ws.Range("C7").Value = "张三"
$ws.Range("D7").Value = "王五"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 32
$ws.Range("G7").Value = 77

$ws.Range("H9").Select()
